$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values per row (F column is unchanged; G = B+C+D+E)
$data = @{
    2  = @(0.003078177322033415, 0.002658071450198252, 0.1496068669990043, 0.5333859586016987)
    3  = @(3.272327238179451,    1.626987699542094,    0.1496068669990043, 0.5333859586016987)
    4  = @(3.272327238179451,    1.626987699542094,    0.7210945179870265, 0.5333859586016987)
    5  = @(3.272327238179451,    1.626987699542094,    0.7210945179870265, 0.5333859586016987)
    6  = @(0.6545652718822623,   1.626987699542094,    0.7210945179870265, 0.5333859586016987)
    7  = @(0.04172184405617529,  0.3048912486333797,   0.7210945179870265, 0.5333859586016987)
    8  = @(3.272327238179451,    1.626987699542094,    0.7210945179870265, 0.5333859586016987)
    9  = @(0.6545652718822623,   0.3048912486333797,   0.7210945179870265, 0.5333859586016987)
    10 = @(3.272327238179451,    1.626987699542094,    0.1496068669990043, 0.5333859586016987)
    11 = @(0.6545652718822623,   0.002658071450198252, 0.7210945179870265, 0.5333859586016987)
    12 = @(0.6545652718822623,   1.626987699542094,    0.7210945179870265, 0.5333859586016987)
    13 = @(3.272327238179451,    1.626987699542094,    0.1496068669990043, 0.5333859586016987)
    14 = @(1.445647641019636,    1.626987699542094,    0.1496068669990043, 0.5333859586016987)
    15 = @(1.445647641019636,    1.626987699542094,    3.223369029078222,  0.5333859586016987)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $b + $c + $d + $e
}
